$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-19 down to 16-20
$ws.Rows.Item(15).Insert()

# Fill in the new row 15 with the new weekly data point
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value = "Maule"
$ws.Cells.Item(15, 4).Value = 44711
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = 100112040
$ws.Cells.Item(15, 7).Value = "Cilantro"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 150
$ws.Cells.Item(15, 11).Value = 8500
$ws.Cells.Item(15, 12).Value = 8500
$ws.Cells.Item(15, 13).Value = 8500
$ws.Cells.Item(15, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 236
$ws.Cells.Item(15, 17).Value = 36
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Apply the same number format (date number format) used on other D column cells
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
